$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 621, shifting existing rows 621-718 down to 623-720.
$ws.Rows("621:622").Insert(-4121)  # -4121 = xlShiftDown

# Populate the two newly inserted rows with the new data.

# Row 621
$ws.Range("A621").Value = 10
$ws.Range("B621").Value = "Vega Modelo de Temuco"
$ws.Range("C621").Value = "La Araucanía"
$ws.Range("D621").Value = 44984
$ws.Range("E621").Value = 9
$ws.Range("F621").Value = 100112043
$ws.Range("G621").Value = "Pepino ensalada"
$ws.Range("H621").Value = "Sin especificar"
$ws.Range("I621").Value = "Primera"
$ws.Range("J621").Value = 370
$ws.Range("K621").Value = 10000
$ws.Range("L621").Value = 11000
$ws.Range("M621").Value = 10595
$ws.Range("N621").Value = "$/caja 60 unidades"
$ws.Range("O621").Value = "Región de Arica y Parinacota"
$ws.Range("P621").Value = 177
$ws.Range("Q621").Value = 60
$ws.Range("R621").Value = "Hortaliza"

# Row 622
$ws.Range("A622").Value = 10
$ws.Range("B622").Value = "Vega Modelo de Temuco"
$ws.Range("C622").Value = "La Araucanía"
$ws.Range("D622").Value = 44984
$ws.Range("E622").Value = 9
$ws.Range("F622").Value = 100112043
$ws.Range("G622").Value = "Pepino ensalada"
$ws.Range("H622").Value = "Sin especificar"
$ws.Range("I622").Value = "Primera"
$ws.Range("J622").Value = 500
$ws.Range("K622").Value = 10000
$ws.Range("L622").Value = 12000
$ws.Range("M622").Value = 10800
$ws.Range("N622").Value = "$/caja 60 unidades"
$ws.Range("O622").Value = "Región del Maule"
$ws.Range("P622").Value = 180
$ws.Range("Q622").Value = 60
$ws.Range("R622").Value = "Hortaliza"

# Ensure the date cells keep the same date number-format style as the rest of column D.
$ws.Range("D621:D622").NumberFormat = $ws.Range("D623").NumberFormat
